$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "м. К" + "иїв, ... ${Contract date}" runs get merged into a single run.
#    The visible text is unchanged (the concatenation already reads
#    "м. Київ,       ${Contract date}"), so replacing the spanning text with
#    itself makes the engine recombine the two runs into one.
# ---------------------------------------------------------------------------
$f1 = $d.Content.Find
$f1.Execute('м. Київ', $false, $false, $false, $false, $false, $true, 1, $false, 'м. Київ', 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) The following run (a single non-breaking space) changes its color from
#    "auto" to 00000A. Locate it precisely as the single character right
#    after "${Contract date}".
# ---------------------------------------------------------------------------
$f2 = $d.Content.Find
$f2.Execute('${Contract date}') | Out-Null
$rngColor = $f2.Parent
$rngColor.Collapse(0)
$rngColor.MoveEnd(1, 1)
$rngColor.Font.Color = 655360

# ---------------------------------------------------------------------------
# Helper: replace the text of a whole paragraph (found via Find) while
# preserving the paragraph's empty <w:rPr/> run-properties element. Doing a
# plain Find/Replace (or Range.Text=) on the run causes the engine to drop
# the empty <w:rPr/>, so instead a fresh paragraph is produced right after
# the original (which keeps <w:rPr/>), the desired text is inserted into it,
# and the original paragraph is deleted afterwards.
# ---------------------------------------------------------------------------
function Replace-ParagraphText($searchText, $newText) {
    $f = $d.Content.Find
    $f.Execute($searchText) | Out-Null
    $rng = $f.Parent
    $para = $rng.Paragraphs(1)
    $paraRange = $para.Range
    $paraRange.InsertParagraphAfter()
    $newPara = $para.Next()
    $newPara.Range.InsertBefore($newText)
    $para.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Requisites block: rewrite the four placeholder paragraphs. Order matters
#    because each new text reuses a placeholder name that a later rule also
#    searches for - doing them top-to-bottom (in original document order)
#    avoids an already-replaced paragraph being replaced again.
# ---------------------------------------------------------------------------
Replace-ParagraphText '${Passport ID}' 'ПІБ як в паспорті = ${Name}'
Replace-ParagraphText '${Passport address}' 'Серіяномер паспорта = ${Passport ID}'
Replace-ParagraphText '${Passport date}' 'ким виданий паспорт = ${Passport address}'
Replace-ParagraphText '${Address}' 'коли виданий паспорт = ${Passport date}'

# ---------------------------------------------------------------------------
# 4) Insert a brand-new paragraph (same TextBody style/spacing) right after
#    the one that used to hold "${Address}", containing the new field.
# ---------------------------------------------------------------------------
$f3 = $d.Content.Find
$f3.Execute('коли виданий паспорт = ${Passport date}') | Out-Null
$para2 = $f3.Parent.Paragraphs(1)
$para2.Range.InsertParagraphAfter()
$newPara2 = $para2.Next()
$newPara2.Range.InsertBefore('Прописка = ${Address}')

# ---------------------------------------------------------------------------
# 5) Normal style: flip overflowPunct from false to true.
# ---------------------------------------------------------------------------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = $true

# ---------------------------------------------------------------------------
# 6) Append eighteen new character styles ListLabel55 .. ListLabel72.
#    55-63 carry a uk-UA language tag, 64-72 do not.
# ---------------------------------------------------------------------------
for ($i = 55; $i -le 72; $i++) {
    $styleId = "ListLabel$i"
    $styleName = "ListLabel $i"
    $ns = $d.Styles.Add($styleId, 2)
    $ns.NameLocal = $styleName
    $ns.QuickStyle = $true
    $ns.Font.NameBi = "OpenSymbol;Arial Unicode MS"
    if ($i -le 63) {
        $ns.Font.LanguageID = "uk-UA"
    }
}
